# Apply quarterly financial updates to the PLTYF sheet.
# Source: "Doing Updates for Financials" commit - numeric restatements
# across the Income Statement / Balance Sheet / Cash Flow sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("H8").Value = 50200
    $ws.Range("I8").Value = 47200
    $ws.Range("J8").Value = 39400
    $ws.Range("H9").Value = 36400
    $ws.Range("I9").Value = 33700
    $ws.Range("J9").Value = 27700
    $ws.Range("H10").Value = 13800
    $ws.Range("I10").Value = 13500
    $ws.Range("J10").Value = 11600
    $ws.Range("H17").Value = 42200
    $ws.Range("I17").Value = 37700
    $ws.Range("J17").Value = 32300
    $ws.Range("I18").Value = 9500
    $ws.Range("E21").Value = 500
    $ws.Range("H21").Value = 15700
    $ws.Range("J21").Value = 9900
    $ws.Range("I26").Value = 7400
    $ws.Range("I27").Value = 7400
    $ws.Range("E29").Value = 18000
    $ws.Range("G29").Value = 74200
    $ws.Range("G33").Value = 74200
    $ws.Range("I33").Value = 20700
    $ws.Range("G35").Value = 74200
    $ws.Range("I35").Value = 20700
    $ws.Range("D41").Value = 66000
    $ws.Range("E41").Value = 50300
    $ws.Range("F41").Value = 56300
    $ws.Range("G41").Value = 61900
    $ws.Range("H41").Value = 206200
    $ws.Range("I41").Value = 69900
    $ws.Range("J41").Value = 70500
    $ws.Range("E43").Value = 18000
    $ws.Range("H43").Value = 49600
    $ws.Range("I43").Value = 46400
    $ws.Range("J43").Value = 38900
    $ws.Range("I44").Value = 9800
    $ws.Range("D46").Value = 66600
    $ws.Range("E46").Value = 72600
    $ws.Range("F46").Value = 57400
    $ws.Range("G46").Value = 79300
    $ws.Range("H46").Value = 272100
    $ws.Range("I46").Value = 131800
    $ws.Range("J46").Value = 126400
    $ws.Range("D48").Value = 14100
    $ws.Range("E48").Value = 24800
    $ws.Range("F48").Value = 22000
    $ws.Range("G48").Value = 19600
    $ws.Range("H48").Value = 38300
    $ws.Range("I48").Value = 39100
    $ws.Range("J48").Value = 41100
    $ws.Range("D54").Value = 82600
    $ws.Range("E54").Value = 99400
    $ws.Range("F54").Value = 81500
    $ws.Range("G54").Value = 101100
    $ws.Range("H54").Value = 314500
    $ws.Range("I54").Value = 175400
    $ws.Range("J54").Value = 172200
    $ws.Range("H57").Value = 166700
    $ws.Range("I57").Value = 31200
    $ws.Range("J57").Value = 31700
    $ws.Range("J58").Value = 5800
    $ws.Range("H59").Value = 10700
    $ws.Range("I59").Value = 9900
    $ws.Range("H60").Value = 180400
    $ws.Range("I60").Value = 45300
    $ws.Range("J60").Value = 48900
    $ws.Range("H66").Value = 180400
    $ws.Range("I66").Value = 45300
    $ws.Range("J66").Value = 48900
    $ws.Range("D72").Value = 77600
    $ws.Range("E72").Value = 95500
    $ws.Range("F72").Value = 77600
    $ws.Range("G72").Value = 97000
    $ws.Range("H72").Value = 131500
    $ws.Range("I72").Value = 126900
    $ws.Range("J72").Value = 120200
    $ws.Range("D76").Value = 80900
    $ws.Range("E76").Value = 98000
    $ws.Range("F76").Value = 80000
    $ws.Range("G76").Value = 99600
    $ws.Range("H76").Value = 134100
    $ws.Range("I76").Value = 130000
    $ws.Range("J76").Value = 123300
    $ws.Range("G81").Value = 74200
    $ws.Range("I81").Value = 20700
    $ws.Range("J83").Value = 2800
    $ws.Range("G89").Value = 9500
    $ws.Range("I89").Value = 14900
    $ws.Range("J89").Value = 10700
    $ws.Range("D94").Value = 31800
    $ws.Range("F94").Value = 13400
    $ws.Range("G94").Value = 88300
    $ws.Range("H94").Value = 127500
    $ws.Range("D96").Value = -19300
    $ws.Range("F96").Value = -19300
    $ws.Range("G96").Value = -104100
    $ws.Range("D100").Value = -19300
    $ws.Range("F100").Value = -19300
    $ws.Range("G100").Value = -105400
    $ws.Range("D102").Value = 15700
    $ws.Range("F102").Value = -5600
    $ws.Range("H102").Value = 136300
    $ws.Range("J102").Value = 9900

# Row 91 (Other Cashflows from Investing Activities), column G (FY2015-12
# period) is restated from a numeric value to "NA".
$ws.Range("G91").Value = "NA"
